$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")
$ws.Activate()

# Insert a new row before existing row 5 (shifts old rows 5..30 down to 6..31).
# Row.Insert() on this runtime copies the formatting of the row above (row 4),
# which already carries the D/E/R/S/T column styles we need for the new row.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 ("3d_Secure" dataset entry).
$ws.Range("A5").Value = "3d_Secure"
$ws.Range("R5").Value = "'4000000000003220"
$ws.Range("S5").Value = "'06/29"
$ws.Range("T5").Value = 123

# Update the view: scroll so column K is the left-most visible column and
# select S8 (matches the author's saved cursor position).
$ws.Range("S8").Select()
